# Applies the "input/ results/ systematic policy_assessment" commit:
#  - appends a trailing "description" row (row 37) to the results sheet that
#    restates the column headers (with a few of them duplicated many times
#    over, exactly as produced by the source tool that generated the sheet)
#  - refreshes the recomputed "Socio-economic capacity" (I) and
#    "Risk to well-being" (J) columns for every province with the
#    newly-recalculated figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Repeat-Text($unit, $times) {
    $result = ""
    for ($i = 0; $i -lt $times; $i++) {
        $result = $result + $unit
    }
    return $result
}

# Copy the formatting of the last data row (36) down into the new row 37 so
# the "description" cell in column A keeps the same bold/bordered style used
# by the rest of column A.
$ws.Range("A36:J36").Copy() | Out-Null
$ws.Range("A37:J37").PasteSpecial(-4122) | Out-Null

$ws.Range("A37").Value2 = "description"
$ws.Range("B37").Value2 = " Average income in the province"
$ws.Range("C37").Value2 = "Population"
$ws.Range("D37").Value2 = " Hazard (protection)"
$ws.Range("E37").Value2 = Repeat-Text "Exposure, poor people" 102
$ws.Range("F37").Value2 = Repeat-Text "Exposure, non-poor people" 101
$ws.Range("G37").Value2 = Repeat-Text "Asset vulnerability (poor people)" 98
$ws.Range("H37").Value2 = Repeat-Text "Asset vulnerability (non-poor people)" 100
$ws.Range("I37").Value2 = Repeat-Text "Socio-economic capacity" 100
$ws.Range("J37").Value2 = Repeat-Text "Risk to well-being" 97

# Recalculated "Socio-economic capacity" (I) / "Risk to well-being" (J)
# values per province (row number -> new value).
$updates = @{
    2  = @(68.808848541394, 1.3763629302022)
    3  = @(76.2904640631491, 3.77788482951005)
    4  = @(53.5117207304491, 1.97505894678772)
    6  = @(130.236668237378, 0.0364618874429158)
    7  = @(59.9029008991516, 0.256597112848222)
    8  = @(174.397746913654, 0.0639838128353085)
    9  = @(95.735157215736, 0.851476656467316)
    10 = @(64.0916265542141, 1.10351012413379)
    11 = @(90.8086768123949, 0.979684057049211)
    12 = @(68.6352090710821, 0.245501129998845)
    13 = @(74.5496973328711, 0.453544554795311)
    14 = @(127.229555316213, 0.0796343191791082)
    15 = @(102.698254536806, 0.10598260899289)
    16 = @(87.4200766538028, 0.615677168570145)
    17 = @(97.6440364315062, 0.486455245476501)
    18 = @(122.13799541939, 0.0433647117145507)
    19 = @(176.55136507901, 0.0661939008427576)
    20 = @(66.5705490549042, 0.824222321878298)
    21 = @(38.5700075266629, 0.853548109109553)
    22 = @(46.596829995927, 3.19139237499955)
    23 = @(70.5164354826221, 0.670031040412605)
    25 = @(92.7205132736894, 0.928616049318641)
    26 = @(116.159677747957, 0.67910742744851)
    27 = @(159.44434752859, 0.458974766668852)
    28 = @(103.237657145986, 0.513751968474989)
    29 = @(97.1599243859768, 0.114750070288743)
    31 = @(56.5965344860607, 1.17319883024291)
    32 = @(44.76954160938, 0.108827338055115)
    33 = @(88.3628584488853, 0.127560301805581)
    34 = @(47.5020909154704, 1.86397769494471)
    35 = @(112.958268860075, 0.553447562293963)
    36 = @(74.9122166017746, 0.49862525090646)
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $ws.Range("I$row").Value2 = $pair[0]
    $ws.Range("J$row").Value2 = $pair[1]
}

Write-Host "Applied description row and refreshed I/J columns"
